$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column E (the "Original Holding" column). This shifts the old
# column F ("Company Master Id*") left into E's place.
$ws.Columns.Item(5).Delete()
